$wb = $excel.ActiveWorkbook
$wsLoads = $wb.Worksheets.Item("INTERNAL_LOADS")
$wsComfort = $wb.Worksheets.Item("INDOOR_COMFORT")

# --- INTERNAL_LOADS: add new "Ev_kW" column (N), mirroring format of column M ---
$srcCol = $wsLoads.Range("M1:M20")
$dstCol = $wsLoads.Range("N1:N20")
$srcCol.Copy($dstCol)

# Header text for the new column
$wsLoads.Range("N1").Value = "Ev_kW"

# Data rows: EV load defaults to 0 (copy already placed 0s, but set explicitly to be safe)
$wsLoads.Range("N2").Value = 0
$wsLoads.Range("N3").Value = 0
$wsLoads.Range("N4").Value = 0
$wsLoads.Range("N5").Value = 0
$wsLoads.Range("N6").Value = 0
$wsLoads.Range("N7").Value = 0
$wsLoads.Range("N8").Value = 0
$wsLoads.Range("N9").Value = 0
$wsLoads.Range("N10").Value = 0
$wsLoads.Range("N11").Value = 0
$wsLoads.Range("N12").Value = 0
$wsLoads.Range("N13").Value = 0
$wsLoads.Range("N14").Value = 0
$wsLoads.Range("N15").Value = 0
$wsLoads.Range("N16").Value = 0
$wsLoads.Range("N17").Value = 0
$wsLoads.Range("N18").Value = 0
$wsLoads.Range("N19").Value = 0
$wsLoads.Range("N20").Value = 0

# --- Make INTERNAL_LOADS the active sheet/tab, with N1 selected ---
$wsLoads.Activate()
$wsLoads.Range("N1").Select()
